$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.260.20'
$ws.Range("E2").Value = '  +2.91%  '
$ws.Range("D3").Value = '1.812.01'
$ws.Range("E3").Value = '  +4.02%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.74'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +2.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3658'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +1.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.90'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07675'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +3.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.140'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +2.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.03'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +2.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.306'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +3.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.532'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +4.73%  '
$ws.Range("D16").Value = '1.818.62'
$ws.Range("E16").Value = '  +4.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '95.08'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +9.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001080'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06532'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +4.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9995'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("E21").Value = '  +3.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.238'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +2.20%  '
$ws.Range("D23").Value = '28.271.01'
$ws.Range("E23").Value = '  +2.86%  '
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.079'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -10.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.57'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +7.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.67'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("D28").Value = '2.022.57'
$ws.Range("E28").Value = '  +4.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.280'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -2.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.81'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.930'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +4.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09180'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.467'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -5.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.97'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +2.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02347'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +2.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.198'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +2.55%  '
$ws.Range("E38").Value = '  +1.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6572'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +2.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06209'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.195'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.122'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +2.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.427'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9988'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.85'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6104'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +3.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.746'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.85'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.017'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +3.34%  '
$ws.Range("E50").Value = '  +2.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07000'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +2.11%  '
